# ---------------------------------------------------------------------------
# Update "苏州-漫展信息.xlsx" (Suzhou comic-convention info) to the output
# regenerated at commit 456a3b4.
#
#   * Sheet "展览"   (exhibitions)  : refresh "want to go" counters (col F)
#   * Sheet "演出"   (performances) : a new gig was scraped
#       -> "苏州·乐队番同人only live Band Set二次元乐队拼盘" (2024-10-13)
#          gets inserted as the new row 3, pushing the existing
#          "Luca Stricagnoli" row down to row 4
#   * Sheet "全部类型" (all types, union of every category) : same counter
#     refresh as "展览" (shifted by one row) *plus* the same new gig
#     inserted as row 23 (pushing everything below it down by one row)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# =============================================================================
# Sheet "展览" (exhibitions) -- worksheet #1
# =============================================================================
$wsExpo = $wb.Worksheets.Item(1)

$wsExpo.Range("F3").Value  = 12962   # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$wsExpo.Range("F5").Value  = 80      # 苏州·明日方舟同人展ONLY
$wsExpo.Range("F7").Value  = 47      # 苏州·首届盗墓笔记同人only 吴邪带我回家（聚会）
$wsExpo.Range("F10").Value = 12927   # 苏州·I COME ACG动漫品牌博览会
$wsExpo.Range("F11").Value = 290     # 苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场
$wsExpo.Range("F12").Value = 42      # 苏州·恋与深空only 同人周边套餐
$wsExpo.Range("F13").Value = 8704    # 苏州·理想乡动漫游戏展-两馆全开+三馆间通道
$wsExpo.Range("F14").Value = 7711    # 苏州·第四届-OCG国朝动漫游戏嘉年华
$wsExpo.Range("F15").Value = 201     # 常熟·CDW·动漫展06
$wsExpo.Range("F16").Value = 114     # 苏州·授渔动漫国风展2.5
$wsExpo.Range("F18").Value = 129     # 苏州·第二届百合Only同人展交流
$wsExpo.Range("F19").Value = 988     # 常熟·cc动漫游戏嘉年华
$wsExpo.Range("F24").Value = 325     # 苏州·绘时国乙2.0光夜同人only

# =============================================================================
# Sheet "演出" (performances) -- worksheet #2
# New gig scraped: inserted above the Luca Stricagnoli row, which shifts
# down from row 3 to row 4 (dimension grows A1:I3 -> A1:I4).
# =============================================================================
$wsShow = $wb.Worksheets.Item(2)

$wsShow.Rows.Item(3).Insert()

# the row that got pushed down keeps its data but its running index (col A)
# advances by one, same as every other row in this sheet
$wsShow.Range("A4").Value = 3

# clone column-A's number/border/alignment formatting onto the fresh row
$wsShow.Range("A2").Copy($wsShow.Range("A3"))

$wsShow.Range("A3").Value = 2
$wsShow.Range("B3").NumberFormat = "@"
$wsShow.Range("B3").Value = "2024-10-13"
$wsShow.Range("C3").Value = "苏州·乐队番同人only live Band Set二次元乐队拼盘"
$wsShow.Range("D3").Value = "扬富路9号南岸新地一期NanNan Park三楼 Wave Livehouse（南岸店）"
$wsShow.Range("E3").Value = "2024.10.13 13:00-10.13 15:30"
$wsShow.Range("F3").Value = 0
$wsShow.Range("G3").Value = 88
$wsShow.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=91594"
$wsShow.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/dPXuP7Q71725003286301.jpeg"

# =============================================================================
# Sheet "全部类型" (all types) -- worksheet #4
# Same counter refresh as "展览" (rows offset by +1 here), then the same new
# gig is inserted above row 23 (the Luca Stricagnoli row), pushing rows
# 23-28 down to 24-29 (dimension grows A1:I28 -> A1:I29).
# =============================================================================
$wsAll = $wb.Worksheets.Item(4)

$wsAll.Range("F4").Value  = 12962   # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$wsAll.Range("F6").Value  = 80      # 苏州·明日方舟同人展ONLY
$wsAll.Range("F8").Value  = 47      # 苏州·首届盗墓笔记同人only 吴邪带我回家（聚会）
$wsAll.Range("F11").Value = 12927   # 苏州·I COME ACG动漫品牌博览会
$wsAll.Range("F12").Value = 290     # 苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场
$wsAll.Range("F13").Value = 42      # 苏州·恋与深空only 同人周边套餐
$wsAll.Range("F14").Value = 8704    # 苏州·理想乡动漫游戏展-两馆全开+三馆间通道
$wsAll.Range("F15").Value = 7711    # 苏州·第四届-OCG国朝动漫游戏嘉年华
$wsAll.Range("F16").Value = 201     # 常熟·CDW·动漫展06
$wsAll.Range("F17").Value = 114     # 苏州·授渔动漫国风展2.5
$wsAll.Range("F19").Value = 129     # 苏州·第二届百合Only同人展交流
$wsAll.Range("F20").Value = 988     # 常熟·cc动漫游戏嘉年华

$wsAll.Rows.Item(23).Insert()

# everything that shifted down (old rows 23..28 -> new rows 24..29) keeps its
# running index one higher than before, same as every other row in this sheet
$wsAll.Range("A24").Value = 23
$wsAll.Range("A25").Value = 24
$wsAll.Range("A26").Value = 25
$wsAll.Range("A27").Value = 26
$wsAll.Range("A28").Value = 27
$wsAll.Range("A29").Value = 28

# the "想去人数" counter for 绘时国乙2.0光夜同人only also refreshed, same as
# on the "展览" sheet above; it now lives on row 27
$wsAll.Range("F27").Value = 325

# clone column-A's number/border/alignment formatting onto the fresh row
$wsAll.Range("A22").Copy($wsAll.Range("A23"))

$wsAll.Range("A23").Value = 22
$wsAll.Range("B23").NumberFormat = "@"
$wsAll.Range("B23").Value = "2024-10-13"
$wsAll.Range("C23").Value = "苏州·乐队番同人only live Band Set二次元乐队拼盘"
$wsAll.Range("D23").Value = "扬富路9号南岸新地一期NanNan Park三楼 Wave Livehouse（南岸店）"
$wsAll.Range("E23").Value = "2024.10.13 13:00-10.13 15:30"
$wsAll.Range("F23").Value = 0
$wsAll.Range("G23").Value = 88
$wsAll.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=91594"
$wsAll.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202408/dPXuP7Q71725003286301.jpeg"
